# Add another response row for the "ทำอะไร" tag.
# Row 28 is currently blank/unused (the sheet's rows are sparse, jumping
# from 27 straight to 29), so we just populate it in place -- no rows
# need to be inserted or shifted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill row 28 with the tag (same as row 27: "ทำอะไร") and the new response.
$ws.Range("A28").Value = "ทำอะไร"
$ws.Range("B28").Value = "ให้ความรู้เรื่องกฎหมายการขายของออนไลน์น้าบ"
